# Charcoal_stoves.xlsx - "First test showing that phi.u is NA."
#
# The FIN_ETA sheet's row 4 describes the "phi.u" (exergy-to-energy ratio)
# quantity for GHA/Final/PCM/phi.u. Column H (the 1971 figure) is cleared
# out entirely so the cell reads as NA, while column I (the 2000 figure)
# keeps its value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FIN_ETA")
$ws.Activate()

# Clear H4 (phi.u, 1971 column) completely -- value AND formatting --
# so it truly drops out of the sheet rather than leaving an empty,
# styled placeholder cell behind.
$ws.Range("H4").Clear()

# The author's selection ended up on H5 after making the edit.
[void]$ws.Range("H5").Select()
